$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header row: the "_old" column group becomes "_FV2210",
# the "_new" column group becomes "_FV2304" (the "diff" column, K1, is unchanged).
$baseNames = @("Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID", "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung")

$oldCols = @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J")
for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Range($oldCols[$i] + "1").Value = $baseNames[$i] + "_FV2210"
}

$newCols = @("L", "M", "N", "O", "P", "Q", "R", "S", "T", "U")
for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Range($newCols[$i] + "1").Value = $baseNames[$i] + "_FV2304"
}

# Turn the data range into an Excel Table (ListObject) with an AutoFilter,
# covering the full used range A1:U57.
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U57"), $null, 1)
$tbl.Name = "Table1"

# Freeze the header row (split after row 1).
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
